$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row data: label, col B, col C, col D, col E
$ws.Range("A10").Value = "Original monthly epoch, old split, random seed = 11"
$ws.Range("B10").Value = 83.19
$ws.Range("C10").Value = 73.37
$ws.Range("D10").Value = 59.83
$ws.Range("E10").Value = 72.6

$ws.Range("A11").Value = "Original monthly epoch, old split, random seed = 77"
$ws.Range("B11").Value = 85.04
$ws.Range("C11").Value = 66.29
$ws.Range("D11").Value = 68.27
$ws.Range("E11").Value = 89.87

$ws.Range("A12").Value = "Original monthly epoch, old split, random seed = 93"
$ws.Range("B12").Value = 85.12
$ws.Range("C12").Value = 69.19
$ws.Range("D12").Value = 76.62
$ws.Range("E12").Value = 83.87

# Copy style from A8 (same label style) to A10:A12
$ws.Range("A8").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122) | Out-Null

# Update selection to match target state
$ws.Range("D29:I41").Select()
$wb.Application.ActiveWindow.RangeSelection.Item(1).Activate()
